# Fruta / hortaliza, semanal
# Insert a new weekly observation at row 71, pushing the existing rows
# 71-95 down to 72-96 (row 96 becomes a new last row, identical to the
# old row 95), and populate the newly inserted row 71 with this week's
# data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 71:95 down to 72:96, inserting a fresh blank row at 71.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new observation.
$ws.Cells.Item(71, 1).Value = 8
$ws.Cells.Item(71, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(71, 3).Value = "Coquimbo"
$ws.Cells.Item(71, 4).Value = 44508
$ws.Cells.Item(71, 5).Value = 4
$ws.Cells.Item(71, 6).Value = 100112001
$ws.Cells.Item(71, 7).Value = "Berenjena"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 520
$ws.Cells.Item(71, 11).Value = 8000
$ws.Cells.Item(71, 12).Value = 9000
$ws.Cells.Item(71, 13).Value = 8500
$ws.Cells.Item(71, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(71, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(71, 16).Value = 142
$ws.Cells.Item(71, 17).Value = 60
$ws.Cells.Item(71, 18).Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of the
# column (style index 2 on D was already inherited by Insert(), but make
# sure explicitly in case the host didn't carry it over).
$ws.Cells.Item(71, 4).NumberFormat = $ws.Cells.Item(72, 4).NumberFormat
